# Feria Lagunitas de Puerto Montt - Uva
# Insert a new weekly price-report row (Red Globe, Provincia de San Felipe
# de Aconcagua) before the current row 173. Excel's native row insert
# shifts every existing row at/after 173 down by one (A1:T255 -> A1:T256)
# and the sheet's dimension / row indices update automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(173).Insert()

$ws.Range("A173").Value = 4
$ws.Range("B173").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C173").Value = "Los Lagos"
$ws.Range("D173").Value = 44806
$ws.Range("E173").Value = 10
$ws.Range("F173").Value = "Fruta"
$ws.Range("G173").Value = 100109
$ws.Range("H173").Value = "Uva"
$ws.Range("I173").Value = 100109001
$ws.Range("J173").Value = "Uva"
$ws.Range("K173").Value = "Red Globe"
$ws.Range("L173").Value = "Primera"
$ws.Range("M173").Value = 300
$ws.Range("N173").Value = 16000
$ws.Range("O173").Value = 17000
$ws.Range("P173").Value = 16500
$ws.Range("Q173").Value = "`$/bandeja 8 kilos"
$ws.Range("R173").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S173").Value = 2062
$ws.Range("T173").Value = 8
